# Added 4wk low sales check - update forecast comparison metrics and summary totals

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Inventory Coverage (H) values per row; $null clears the cell
$hValues = @{
    2  = 9
    3  = 8
    4  = 7
    5  = 6
    6  = $null
    7  = $null
    8  = $null
    9  = $null
    10 = $null
    11 = $null
    12 = 5
    13 = 4
    14 = $null
    15 = $null
    16 = $null
    17 = $null
}

# New MyForecast (D) values per row (only rows that change)
$dValues = @{
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
}

# New Stockout Risk (I) and Reorder Urgency (J) values per row (only rows that change)
$iValues = @{
    7  = "Low"
    8  = "Low"
    9  = "Low"
    10 = "Low"
    11 = "Low"
    12 = "Low"
    13 = "Low"
    14 = "Low"
    15 = "Low"
    16 = "Low"
    17 = "Low"
}
$jValues = @{
    7  = "Normal"
    8  = "Normal"
    9  = "Normal"
    10 = "Normal"
    11 = "Normal"
    12 = "Normal"
    13 = "Normal"
    14 = "Normal"
    15 = "Normal"
    16 = "Normal"
    17 = "Normal"
}

# New Seasonality Index (L) values per row
$lValues = @{
    2  = 1.04
    3  = 1.12
    4  = 1.15
    5  = 0.99
    6  = 1
    7  = 1.09
    8  = 1.2
    9  = 0.85
    10 = 1.12
    11 = 1.15
    12 = 1.15
    13 = 1.19
    14 = 0.87
    15 = 0.84
    16 = 1.06
    17 = 0.89
}

foreach ($row in $dValues.Keys) {
    $ws1.Cells.Item($row, 4).Value = $dValues[$row]
}

foreach ($row in $hValues.Keys) {
    if ($null -eq $hValues[$row]) {
        $ws1.Cells.Item($row, 8).ClearContents()
    } else {
        $ws1.Cells.Item($row, 8).Value = $hValues[$row]
    }
}

foreach ($row in $iValues.Keys) {
    $ws1.Cells.Item($row, 9).Value = $iValues[$row]
}

foreach ($row in $jValues.Keys) {
    $ws1.Cells.Item($row, 10).Value = $jValues[$row]
}

foreach ($row in $lValues.Keys) {
    $ws1.Cells.Item($row, 12).Value = $lValues[$row]
}

# --- Sheet 2: "Summary" ---
# These values are stored as text (not numbers) in the workbook, so a leading
# apostrophe is used to force Excel to keep them as text, matching the
# original inline-string cell type.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value  = "'6"
$ws2.Range("B10").Value = "'4"
$ws2.Range("B11").Value = "'4"
$ws2.Range("B12").Value = "'1"
$ws2.Range("B14").Value = "'0"
